# Add 2022-Q4 data.
#
# The workbook currently has three sheets: "总计" (summary), "2022-Q2" and
# "2021-Q1". We add a new "2022-Q4" sheet (with its own fund-holdings
# table) positioned right after "总计" and before "2022-Q2", and record it
# in the "总计" summary sheet. The previously-existing "2022-Q2" sheet and
# its data are preserved unchanged, just shifted one position to the right.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: duplicate the existing "2022-Q2" sheet, placing the copy right
# after it. The copy keeps the old 2022-Q2 figures untouched; the
# original slot is turned into the new "2022-Q4" sheet below. This keeps
# the "2022-Q2" data/formatting completely intact while freeing up the
# original sheet object to become the newly reported quarter.
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($null, $q2)
$q2Dup = $wb.Worksheets.Item("2022-Q2 (2)")

# ---------------------------------------------------------------------
# Step 2: rewrite the original "2022-Q2" sheet object in place with the
# new 2022-Q4 fund holdings, then rename it. Values are entered with a
# leading "'" so fund codes like "005189" and ratios like "1.00" are
# kept as literal text (matching the source data) instead of being
# auto-coerced into numbers; ClearFormats() afterwards drops the
# resulting quote-prefix styling so the cells stay in their original,
# un-styled state.
# ---------------------------------------------------------------------
$q4 = $q2
$q4.Cells.Item(2, 2).Value = "'005189"
$q4.Cells.Item(2, 3).Value = "'海富通量化前锋股票A"
$q4.Cells.Item(2, 4).Value = "'0.54"
$q4.Cells.Item(2, 5).Value = "'88.13"
$q4.Cells.Item(2, 6).Value = "'1.00"
$q4.Cells.Item(2, 7).Value = "'0.0054"
$q4.Cells.Item(2, 8).Value = 8

$q4.Cells.Item(3, 2).Value = "'005188"
$q4.Cells.Item(3, 3).Value = "'海富通量化前锋股票C"
$q4.Cells.Item(3, 4).Value = "'0.03"
$q4.Cells.Item(3, 5).Value = "'88.13"
$q4.Cells.Item(3, 6).Value = "'1.00"
$q4.Cells.Item(3, 7).Value = "'0.0003"
$q4.Cells.Item(3, 8).Value = 8

$q4.Range("B2:G3").ClearFormats()

$q4.Name = "2022-Q4"

# Restore the duplicate's name back to "2022-Q2" - it still holds the
# original, unmodified 2022-Q2 figures.
$q2Dup.Name = "2022-Q2"

# ---------------------------------------------------------------------
# Step 3: update the "总计" overview sheet so it lists 2022-Q4 first,
# then the (now shifted down) 2022-Q2 and 2021-Q1 rows.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Existing row 2 ("2022-Q2") becomes the new 2022-Q4 summary row.
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 2
$total.Cells.Item(2, 4).Value = 0.01

# Copy row 3's formatting down into a new row 4 before repurposing row 3,
# so the inserted "2021-Q1" row keeps the same cell styling (A4 etc.)
$total.Range("A3:D3").Copy()
$total.Range("A4:D4").PasteSpecial(-4122)

# Row 4 becomes the (shifted) "2021-Q1" entry.
$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(4, 2).Value = "2021-Q1"
$total.Cells.Item(4, 3).Value = 1
$total.Cells.Item(4, 4).Value = 0.69

# Row 3 becomes the (shifted) "2022-Q2" entry.
$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(3, 2).Value = "2022-Q2"
$total.Cells.Item(3, 3).Value = 2
$total.Cells.Item(3, 4).Value = 0.07

# "2021-Q1" was (and remains) the last/selected tab - restore that after
# all the sheet shuffling above.
$wb.Worksheets.Item("2021-Q1").Activate()
